# Adds a new "2020" column (Q) to the data table, mirroring the
# formatting of the existing "2019" column (P), and updates the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 header: new year value in Q3
$Q3 = 2020

# Data rows 4-33: new values for column Q, in row order
$values = [ordered]@{
    4  = 1.9148453093736542
    5  = 1.7453236044300597
    6  = 2.0818900906859255
    7  = 1.658050942694075
    8  = 1.4467487937731931
    9  = 1.8774124750304142
    10 = 0.96024351775610284
    11 = 0.63595936855594293
    12 = 1.2888424905592288
    13 = 1.6032353288937073
    14 = 2.4146715443031859
    15 = 0.79837132250209564
    16 = 1.3751327862596732
    17 = 0.67516929870164943
    18 = 2.1012817818869509
    19 = 1.5943738893736428
    20 = 1.5765365498500856
    21 = 1.6126194804433236
    22 = 0.37150276583809166
    23 = 0
    24 = 0.75125835774923
    25 = 2.8942542850468351
    26 = 2.72898263527357
    27 = 3.0545792215303034
    28 = 3.9473869708034344
    29 = 3.6031203021816895
    30 = 4.2520923837938582
    31 = 0
    32 = 0
    33 = 0
}

# Copy formatting from column P to column Q for rows 3-33 so the new
# column inherits the same number format / borders / font as the
# "2019" column it follows.
$srcRange = $ws.Range("P3:P33")
$dstRange = $ws.Range("Q3:Q33")
$srcRange.Copy() | Out-Null
$dstRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Write the header value
$ws.Range("Q3").Value = $Q3

# Write the data values
foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 17).Value = $values[$row]
}

# Update the active cell / selection to T1, matching the post-edit
# worksheet state.
$ws.Range("T1").Select() | Out-Null
